$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Apply the same style (centered) used by the existing data rows to the new rows A:C ---
$ws.Range("A9:C11").HorizontalAlignment = -4108

# --- Row 9: Cambridge University Press Full Package ---
$ws.Cells.Item(9, 1).Value = 8
$ws.Cells.Item(9, 2).Value = 1
$ws.Cells.Item(9, 3).Value = "kbart"
$ws.Cells.Item(9, 4).Value = "CUP_FP"
$ws.Cells.Item(9, 5).Value = "Cambridge Journals Full Package"
$ws.Cells.Item(9, 6).Value = "https://www.cambridge.org/core/services/aop-cambridge-core/kbart/create/bespoke/BE6D264D98C2E9F9BFDC17C422C3C696"

# --- Row 10: Oxford University Press Current Collection ---
$ws.Cells.Item(10, 1).Value = 9
$ws.Cells.Item(10, 2).Value = 1
$ws.Cells.Item(10, 3).Value = "kbart"
$ws.Cells.Item(10, 4).Value = "OUP_CC"
$ws.Cells.Item(10, 5).Value = "Oxford Current Collection"
$ws.Cells.Item(10, 6).Value = "http://fdslive.oup.com/www.oup.com/academic/content/librarian/OxfordUniversityPress_Global_2021JournalsCurrentCollection.zip"

# --- Row 11: Oxford University Press Open Access Tiles ---
$ws.Cells.Item(11, 1).Value = 10
$ws.Cells.Item(11, 2).Value = 1
$ws.Cells.Item(11, 3).Value = "kbart"
$ws.Cells.Item(11, 4).Value = "OUP_OA"
$ws.Cells.Item(11, 5).Value = "Oxford Open Access Tiles"
$ws.Cells.Item(11, 6).Value = "http://fdslive.oup.com/www.oup.com/academic/content/librarian/OxfordUniversityPress_Global_2021JournalsOpenAccess.zip"

# --- Update existing rows: Springer/Nature codes get an "SN_" prefix ---
$ws.Cells.Item(5, 4).Value = "SN_SPRINGER"
$ws.Cells.Item(6, 4).Value = "SN_NATURE"

# --- Update selection / view state ---
$null = $ws.Range("D20").Select()
$excel.ActiveWindow.Left = 1140
$excel.ActiveWindow.Top = 2100
